# Update cryptocurrency price (column D) and 1h volume change (column E)
# values on Sheet1 to reflect the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> @{ D = "<new price>"; E = "<new change>" }
# Only keys present are updated (some rows only change D or only E).
$updates = @{
    2  = @{ D = "273.39";     E = "-1.81%" }
    3  = @{ D = "26.65";      E = "-2.35%" }
    4  = @{ D = "4.886";      E = "1.96%"  }
    5  = @{             E = "1.38%"  }
    6  = @{ D = "6.902";      E = "0.83%"  }
    7  = @{ D = "3.357";      E = "2.58%"  }
    8  = @{ D = "1.265";      E = "34.95%" }
    9  = @{ D = "0.8777";     E = "-0.01%" }
    10 = @{ D = "0.1452";     E = "-0.07%" }
    11 = @{ D = "0.05013";    E = "-0.59%" }
    12 = @{ D = "0.07324";    E = "0.56%"  }
    13 = @{ D = "0.03119";    E = "-0.98%" }
    14 = @{ D = "0.09048";    E = "0.14%"  }
    15 = @{ D = "0.001597";   E = "3.17%"  }
    16 = @{ D = "0.0006293";  E = "0.44%"  }
    17 = @{ D = "0.006070";   E = "-0.65%" }
    18 = @{             E = "-0.43%" }
    19 = @{             E = "0.74%"  }
    20 = @{             E = "2.56%"  }
    21 = @{ D = "0.1326";     E = "1.22%"  }
    22 = @{ D = "3.919";      E = "1.75%"  }
    23 = @{             E = "1.99%"  }
    24 = @{ D = "0.001180";   E = "0.33%"  }
    25 = @{ D = "0.004405";   E = "3.24%"  }
    26 = @{ D = "0.0001202";  E = "0.28%"  }
    27 = @{ D = "0.0001703";  E = "5.54%"  }
    40 = @{ D = "0.04028";    E = "0.03%"  }
    41 = @{ D = "0.006618";   E = "-1.23%" }
    42 = @{ D = "0.1166";     E = "1.39%"  }
    43 = @{ D = "0.002104";   E = "-1.12%" }
    44 = @{ D = "0.01253";    E = "1.59%"  }
    45 = @{ D = "0.00005316"; E = "3.80%"  }
    46 = @{             E = "11.12%" }
    47 = @{ D = "0.02002";    E = "-32.97%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").NumberFormat = "@"
        $ws.Range("E$row").Value = $vals["E"]
    }
}
